# terneiraTrack workbook update:
#  - populate "2- crescimento" (growth) sheet with calf growth-tracking data
#  - add a new "3- saude" (health) sheet with a health-tracking header row
#  - make the new sheet the active tab

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet 2 ("2- crescimento") - growth tracking table
# ---------------------------------------------------------------------

# Header row
$ws2.Range("A1").Value = "Num. Bezerra"
$ws2.Range("B1").Value = "Data nasc."
$ws2.Range("C1").Value = "Peso Nasc. (kg)"
$ws2.Range("D1").Value = "Data Pesagem"
$ws2.Range("E1").Value = "Peso (kg)"

# Column A - calf number
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
$ws2.Range("A6").Value = 5

# Column B - birth date (date-formatted); format the whole block B2:B7
# (including the trailing blank row) via a single copy/paste so every
# cell shares one style entry, matching a hand-formatted column.
$ws2.Range("B2").Value = 45394
$ws2.Range("B2").NumberFormat = "mm-dd-yy"
$ws2.Range("B2").Copy()
$ws2.Range("B3:B7").PasteSpecial(-4122)
$ws2.Range("B3").Value = 45395
$ws2.Range("B4").Value = 45396
$ws2.Range("B5").Value = 45397
$ws2.Range("B6").Value = 45398

# Column C - birth weight (kg)
$ws2.Range("C2").Value = 35
$ws2.Range("C3").Value = 32
$ws2.Range("C4").Value = 40
$ws2.Range("C5").Value = 42
$ws2.Range("C6").Value = 40

# Column D - weigh-in date = birth date + 30 days (date-formatted formula)
$ws2.Range("D2:D6").PasteSpecial(-4122)
$ws2.Range("D2").Formula = "=B2+30"
$ws2.Range("D3").Formula = "=B3+30"
$ws2.Range("D4").Formula = "=B4+30"
$ws2.Range("D5").Formula = "=B5+30"
$ws2.Range("D6").Formula = "=B6+30"

# Column E - weight at weigh-in = birth weight + 30*0.5
$ws2.Range("E2").Formula = "=C2+30*0.5"
$ws2.Range("E3").Formula = "=C3+30*0.5"
$ws2.Range("E4").Formula = "=C4+30*0.5"
$ws2.Range("E5").Formula = "=C5+30*0.5"
$ws2.Range("E6").Formula = "=C6+30*0.5"

# Column widths (auto-fit look)
$ws2.Columns.Item(2).ColumnWidth = 9
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws2.Columns.Item(4).ColumnWidth = 11.8

# ---------------------------------------------------------------------
# Sheet 3 ("3- saude") - new health tracking sheet, placed after sheet 2
# ---------------------------------------------------------------------

# Clone "1- cadastro" (it carries no stray <cols> overrides) so the new
# sheet inherits the workbook's normal sheetFormatPr (16pt default row
# height) instead of a blank-sheet default, then wipe its contents.
$ws1.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Clear()
$ws3.Name = "3- saúde"

# Assign in this order so new shared-string entries land in the same
# sequence as the authored workbook (Num. Bezerra is already shared
# with sheet 2, so it doesn't allocate a new entry here).
$ws3.Range("A1").Value = "Num. Bezerra"
$ws3.Range("C1").Value = "Data Tratamento"
$ws3.Range("D1").Value = "Tratmento"
$ws3.Range("B1").Value = "Doença"
$ws3.Range("E1").Value = "Obs."

$ws3.Columns.Item(1).ColumnWidth = 11.3
$ws3.Columns.Item(2).ColumnWidth = 11.3
$ws3.Columns.Item(3).ColumnWidth = 14.0

# Leave the cursor parked at E7 on the new sheet, matching the saved
# selection state, and make this the active/visible tab.
$ws3.Activate()
$ws3.Range("E7").Select()
